$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2449.875
$ws.Range("I131").Value = 2014.1428
$ws.Range("K131").Value = 6042.428400000001
$ws.Range("M131").Value = -1002.428400000001
$ws.Range("H141").Value = 1271.3334
$ws.Range("I141").Value = 1279.625
$ws.Range("J141").Value = 1205
$ws.Range("K141").Value = 3838.875
$ws.Range("L141").Value = 3615
$ws.Range("M141").Value = 1341.125
$ws.Range("N141").Value = -13975

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2440
$ws.Range("I61").Value = 2222.5334
$ws.Range("K61").Value = 2222.5334
$ws.Range("M61").Value = -2010.5334
$ws.Range("H74").Value = 2743.6206
$ws.Range("I74").Value = 1403.3077
$ws.Range("J74").Value = 3832.625
$ws.Range("K74").Value = 1403.3077
$ws.Range("L74").Value = 3832.625
$ws.Range("M74").Value = -529.3077000000001
$ws.Range("N74").Value = -5580.625
$ws.Range("H77").Value = 2743.6206
$ws.Range("I77").Value = 1403.3077
$ws.Range("J77").Value = 3832.625
$ws.Range("K77").Value = 7016.538500000001
$ws.Range("L77").Value = 19163.125
$ws.Range("M77").Value = -2648.538500000001
$ws.Range("N77").Value = -27899.125
$ws.Range("H132").Value = 2131.456
$ws.Range("I132").Value = 1937.76
$ws.Range("J132").Value = 3515
$ws.Range("K132").Value = 5813.28
$ws.Range("L132").Value = 10545
$ws.Range("M132").Value = -3283.28
$ws.Range("N132").Value = -15605
$ws.Range("H136").Value = 2440
$ws.Range("I136").Value = 2222.5334
$ws.Range("K136").Value = 6667.600199999999
$ws.Range("M136").Value = -4117.600199999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3132.625
$ws.Range("J86").Value = 3689.25
$ws.Range("L86").Value = 3689.25
$ws.Range("N86").Value = -5935.25
$ws.Range("H89").Value = 3132.625
$ws.Range("J89").Value = 3689.25
$ws.Range("L89").Value = 18446.25
$ws.Range("N89").Value = -29678.25
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H134").Value = 3192
$ws.Range("I134").Value = 1581.4546
$ws.Range("K134").Value = 4744.3638
$ws.Range("M134").Value = -2209.3638

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1681.125
$ws.Range("I16").Value = 1442.2
$ws.Range("K16").Value = 1442.2
$ws.Range("M16").Value = -1155.2
$ws.Range("H58").Value = 1717.56
$ws.Range("I58").Value = 1043.6757
$ws.Range("J58").Value = 3635.5386
$ws.Range("K58").Value = 1043.6757
$ws.Range("L58").Value = 3635.5386
$ws.Range("M58").Value = -840.6757
$ws.Range("N58").Value = -4041.5386
$ws.Range("H113").Value = 1681.125
$ws.Range("I113").Value = 1442.2
$ws.Range("K113").Value = 1442.2
$ws.Range("M113").Value = 727.8
$ws.Range("H132").Value = 50001580
$ws.Range("I132").Value = 1976
$ws.Range("K132").Value = 5928
$ws.Range("M132").Value = -3398
$ws.Range("H136").Value = 1717.56
$ws.Range("I136").Value = 1043.6757
$ws.Range("J136").Value = 3635.5386
$ws.Range("K136").Value = 3131.0271
$ws.Range("L136").Value = 10906.6158
$ws.Range("M136").Value = -581.0271000000002
$ws.Range("N136").Value = -16006.6158

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 237.375
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 237.375
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 712.125
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -992.125
$ws.Range("H99").Value = 6126
$ws.Range("I99").Value = 2675.3333
$ws.Range("K99").Value = 8025.999899999999
$ws.Range("M99").Value = -5779.999899999999
$ws.Range("H108").Value = 9255.25
$ws.Range("I108").Value = 307.33334
$ws.Range("J108").Value = 14624
$ws.Range("K108").Value = 922.0000200000001
$ws.Range("L108").Value = 43872
$ws.Range("M108").Value = 1957.99998
$ws.Range("N108").Value = -49632
$ws.Range("H121").Value = 25001070
$ws.Range("J121").Value = 1362.3334
$ws.Range("L121").Value = 4087.0002
$ws.Range("N121").Value = -6707.0002
$ws.Range("H129").Value = 1037.591
$ws.Range("J129").Value = 1782.8572
$ws.Range("L129").Value = 5348.571599999999
$ws.Range("N129").Value = -15348.5716
$ws.Range("H132").Value = 5049.375
$ws.Range("J132").Value = 5799
$ws.Range("L132").Value = 52191
$ws.Range("N132").Value = -57251

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19794
$ws.Range("I33").Value = 18484.5
$ws.Range("K33").Value = 18484.5
$ws.Range("M33").Value = -18232.5
$ws.Range("H52").Value = 48998.715
$ws.Range("J52").Value = 48998.715
$ws.Range("L52").Value = 48998.715
$ws.Range("N52").Value = -49516.715
$ws.Range("H80").Value = 213396.38
$ws.Range("J80").Value = 4618.6
$ws.Range("L80").Value = 4618.6
$ws.Range("N80").Value = -6614.6
$ws.Range("H83").Value = 213396.38
$ws.Range("J83").Value = 4618.6
$ws.Range("L83").Value = 23093
$ws.Range("N83").Value = -33077
$ws.Range("H126").Value = 3004.1667
$ws.Range("I126").Value = 2936.3635
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 8809.0905
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -6339.0905
$ws.Range("N126").Value = -16190
$ws.Range("H132").Value = 52651900
$ws.Range("I132").Value = 62514132
$ws.Range("J132").Value = 53342.668
$ws.Range("K132").Value = 187542396
$ws.Range("L132").Value = 160028.004
$ws.Range("M132").Value = -187539866
$ws.Range("N132").Value = -165088.004

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 577.38464
$ws.Range("I16").Value = 522.5484
$ws.Range("J16").Value = 789.875
$ws.Range("K16").Value = 522.5484
$ws.Range("L16").Value = 789.875
$ws.Range("M16").Value = -352.5484
$ws.Range("N16").Value = -1129.875
$ws.Range("H22").Value = 1309.2
$ws.Range("I22").Value = 849.3333
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 849.3333
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -554.3333
$ws.Range("N22").Value = -2589
$ws.Range("H27").Value = 1309.2
$ws.Range("I27").Value = 849.3333
$ws.Range("J27").Value = 1999
$ws.Range("K27").Value = 849.3333
$ws.Range("L27").Value = 1999
$ws.Range("M27").Value = -742.3333
$ws.Range("N27").Value = -2213
$ws.Range("H40").Value = 2610.3845
$ws.Range("I40").Value = 2000.7646
$ws.Range("J40").Value = 3761.889
$ws.Range("K40").Value = 2000.7646
$ws.Range("L40").Value = 3761.889
$ws.Range("M40").Value = -1864.7646
$ws.Range("N40").Value = -4033.889
$ws.Range("H132").Value = 2748.913
$ws.Range("I132").Value = 2692.0454
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8076.1362
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5546.1362
$ws.Range("N132").Value = -17060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H132").Value = 2369.8462
$ws.Range("I132").Value = 2064.64
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6193.92
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3663.92
$ws.Range("N132").Value = -35060
